$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 45786.01041666666, 571, 750, 1424.083078233393, 1321),
    @(3, 45786.02083333334, 581, 748, 1423.220003733747, 1329),
    @(4, 45786.03125, 574, 748, 1422.356929234102, 1322),
    @(5, 45786.04166666666, 556, 748, 1421.493854734456, 1304),
    @(6, 45786.05208333334, 551, 756, 1415.574179048649, 1307),
    @(7, 45786.0625, 546, 755, 1409.654503362842, 1301),
    @(8, 45786.07291666666, 537, 755, 1403.734827677036, 1292),
    @(9, 45786.08333333334, 534, 757, 1397.815151991229, 1291),
    @(10, 45786.09375, 548, 757, 1399.246840278499, 1305),
    @(11, 45786.10416666666, 548, 756, 1400.678528565769, 1304),
    @(12, 45786.11458333334, 589, 756, 1402.110216853039, 1345),
    @(13, 45786.125, 574, 756, 1403.541905140309, 1330),
    @(14, 45786.13541666666, 555, 758, 1406.943434051056, 1313),
    @(15, 45786.14583333334, 556, 760, 1410.344962961803, 1316),
    @(16, 45786.15625, 548, 760, 1413.74649187255, 1308),
    @(17, 45786.16666666666, 553, 760, 1417.148020783297, 1313),
    @(18, 45786.17708333334, 550, 762, 1425.839688685614, 1312),
    @(19, 45786.1875, 501, 763, 1434.53135658793, 1264),
    @(20, 45786.19791666666, 540, 765, 1443.223024490247, 1305),
    @(21, 45786.20833333334, 561, 766, 1451.914692392563, 1327),
    @(22, 45786.21875, 341, 798, 1522.382186834008, 1139),
    @(23, 45786.22916666666, 341, 800, 1592.849681275452, 1141),
    @(24, 45786.23958333334, 345, 771, 1663.317175716897, 1116),
    @(25, 45786.25, 383, 760, 1733.784670158341, 1143),
    @(26, 45786.26041666666, 536, 1050, 1749.817548216582, 1586),
    @(27, 45786.27083333334, 0, 0, 1765.850426274822, 0),
    @(28, 45786.28125, 0, 0, 1781.883304333062, 0),
    @(29, 45786.29166666666, 0, 0, 1797.916182391303, 0),
    @(30, 45786.30208333334, 0, 0, 1790.402357335445, 0),
    @(31, 45786.3125, 0, 0, 1782.888532279587, 0),
    @(32, 45786.32291666666, 0, 0, 1775.374707223729, 0),
    @(33, 45786.33333333334, 0, 0, 1767.860882167872, 0),
    @(34, 45786.34375, 0, 0, 1740.547112708727, 0),
    @(35, 45786.35416666666, 0, 0, 1713.233343249583, 0),
    @(36, 45786.36458333334, 0, 0, 1685.919573790439, 0),
    @(37, 45786.375, 0, 0, 1658.605804331294, 0),
    @(38, 45786.38541666666, 0, 0, 1634.175719082284, 0),
    @(39, 45786.39583333334, 0, 0, 1609.745633833273, 0),
    @(40, 45786.40625, 0, 0, 1585.315548584262, 0),
    @(41, 45786.41666666666, 0, 0, 1560.885463335251, 0),
    @(42, 45786.42708333334, 0, 0, 1542.730437389968, 0),
    @(43, 45786.4375, 0, 0, 1524.575411444685, 0),
    @(44, 45786.44791666666, 0, 0, 1506.420385499403, 0),
    @(45, 45786.45833333334, 0, 0, 1488.26535955412, 0),
    @(46, 45786.46875, 0, 0, 1475.644164223811, 0),
    @(47, 45786.47916666666, 0, 0, 1463.022968893501, 0),
    @(48, 45786.48958333334, 0, 0, 1450.401773563191, 0),
    @(49, 45786.5, 0, 0, 1437.780578232882, 0),
    @(50, 45786.51041666666, 0, 0, 1435.699045616186, 0),
    @(51, 45786.52083333334, 0, 0, 1433.61751299949, 0),
    @(52, 45786.53125, 0, 0, 1431.535980382795, 0),
    @(53, 45786.54166666666, 0, 0, 1429.454447766099, 0),
    @(54, 45786.55208333334, 0, 0, 1440.603339537893, 0),
    @(55, 45786.5625, 0, 0, 1451.752231309688, 0),
    @(56, 45786.57291666666, 0, 0, 1462.901123081482, 0),
    @(57, 45786.58333333334, 0, 0, 1474.050014853276, 0),
    @(58, 45786.59375, 0, 0, 1470.638332125364, 0),
    @(59, 45786.60416666666, 0, 0, 1467.226649397453, 0),
    @(60, 45786.61458333334, 0, 0, 1463.814966669541, 0),
    @(61, 45786.625, 0, 0, 1460.403283941629, 0),
    @(62, 45786.63541666666, 0, 0, 1502.328396989123, 0),
    @(63, 45786.64583333334, 0, 0, 1544.253510036617, 0),
    @(64, 45786.65625, 0, 0, 1586.178623084111, 0),
    @(65, 45786.66666666666, 0, 0, 1628.103736131605, 0),
    @(66, 45786.67708333334, 0, 0, 1714.644723901826, 0),
    @(67, 45786.6875, 0, 0, 1801.185711672048, 0),
    @(68, 45786.69791666666, 0, 0, 1887.726699442269, 0),
    @(69, 45786.70833333334, 0, 0, 1974.26768721249, 0),
    @(70, 45786.71875, 0, 0, 2033.383213529679, 0),
    @(71, 45786.72916666666, 0, 0, 2092.498739846868, 0),
    @(72, 45786.73958333334, 0, 0, 2151.614266164058, 0),
    @(73, 45786.75, 0, 0, 2210.729792481247, 0),
    @(74, 45786.76041666666, 0, 0, 2235.606645706179, 0),
    @(75, 45786.77083333334, 0, 0, 2260.483498931112, 0),
    @(76, 45786.78125, 0, 0, 2285.360352156044, 0),
    @(77, 45786.79166666666, 0, 0, 2310.237205380977, 0),
    @(78, 45786.80208333334, 0, 0, 2345.816182401499, 0),
    @(79, 45786.8125, 0, 0, 2381.395159422022, 0),
    @(80, 45786.82291666666, 0, 0, 2416.974136442544, 0),
    @(81, 45786.83333333334, 0, 0, 2452.553113463067, 0),
    @(82, 45786.84375, 0, 0, 2386.868067131375, 0),
    @(83, 45786.85416666666, 0, 0, 2321.183020799682, 0),
    @(84, 45786.86458333334, 0, 0, 2255.49797446799, 0),
    @(85, 45786.875, 0, 0, 2189.812928136298, 0),
    @(86, 45786.88541666666, 0, 0, 2108.237157193317, 0),
    @(87, 45786.89583333334, 0, 0, 2026.661386250336, 0),
    @(88, 45786.90625, 0, 0, 1945.085615307355, 0),
    @(89, 45786.91666666666, 0, 0, 1863.509844364374, 0),
    @(90, 45786.92708333334, 0, 0, 1793.123580463859, 0),
    @(91, 45786.9375, 0, 0, 1722.737316563343, 0),
    @(92, 45786.94791666666, 0, 0, 1652.351052662827, 0),
    @(93, 45786.95833333334, 0, 0, 1581.964788762312, 0),
    @(94, 45786.96875, 0, 0, 1522.425245311292, 0),
    @(95, 45786.97916666666, 0, 0, 1462.885701860272, 0),
    @(96, 45786.98958333334, 0, 0, 1403.346158409252, 0),
    @(97, 45787, 0, 0, 1343.806614958232, 0)
)

foreach ($row in $data) {
    $r = [int]$row[0]
    $ws.Cells.Item($r, 1).Value = [double]$row[1]
    $ws.Cells.Item($r, 2).Value = [double]$row[2]
    $ws.Cells.Item($r, 3).Value = [double]$row[3]
    $ws.Cells.Item($r, 4).Value = [double]$row[4]
    $ws.Cells.Item($r, 5).Value = [double]$row[5]
}
